# Applies the "Responses.docx" review-response edits described by the
# commit: beef up the conclusion response, add a keywords response, and
# merge the word-wrapped runs in the Reviewer 2 paragraph (an artifact of
# that paragraph being retyped).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Response" paragraph right after the "conclusion
#    section" reviewer comment (currently paragraph 7), carrying the
#    _GoBack bookmark that used to sit on paragraph 1.
# ---------------------------------------------------------------------
$pConclusion = $d.Paragraphs.Item(7)
Write-Host "p7 (expect conclusion comment): " $pConclusion.Range.Text.Substring(0, 40)

$pConclusion.Range.InsertParagraphAfter()
$pConclusionResponse = $d.Paragraphs.Item(8)
$pConclusionResponse.Range.Text = "The conclusion is increased 670% and now contains wildly accurate predictions about the future."
$pConclusionResponse.Style = "Response"

# Bookmark goes at the very end of that new paragraph's text (before the
# paragraph mark).
$bmRange = $d.Range($pConclusionResponse.Range.End - 1, $pConclusionResponse.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 2) Insert a new "Response" paragraph right after the "keywords" reviewer
#    comment (paragraph 3).
# ---------------------------------------------------------------------
$pKeywords = $d.Paragraphs.Item(3)
Write-Host "p3 (expect keywords comment): " $pKeywords.Range.Text.Substring(0, 40)

$pKeywords.Range.InsertParagraphAfter()
$pKeywordsResponse = $d.Paragraphs.Item(4)
$pKeywordsResponse.Range.Text = "I" + [char]0x2019 + "ve added keywords for most of the concepts introduced in the paper.  Hopefully I have not gone overboard by introducing too many keywords."
$pKeywordsResponse.Style = "Response"

# ---------------------------------------------------------------------
# 3) Remove the _GoBack bookmark from paragraph 1 ("Reviewer: 1") - it
#    has moved to the new conclusion response paragraph above.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
Write-Host "p1 (expect Reviewer: 1): " $p1.Range.Text

foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack" -and $bm.Start -lt $p1.Range.End) {
        $bm.Delete()
    }
}

# ---------------------------------------------------------------------
# 4) Reviewer 2's "That said, ..." paragraph was retyped, so its many
#    word-wrap-sized runs collapse into two larger runs. Re-create that
#    by deleting everything after "That said," and re-inserting the
#    merged text.
# ---------------------------------------------------------------------
$paragraphs = $d.Paragraphs
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $cand = $paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t.StartsWith("That said")) {
        $pThatSaid = $cand
        break
    }
}

$pStart = $pThatSaid.Range.Start
$paraEnd = $pThatSaid.Range.End
$afterComma = $pStart + 10

$mergedText = " it would be nice if the shortcomings of the pipeline approach were explicitly discussed.  One problem revealed by the survey itself is that in actuality the simple pipeline model is inadequate to describe the recent advances in visualization systems. This is true even in the leading open-source toolkits that implement most of the features discussed (like the VTK). There seems to be no standard notation or pictorial representation of the modern visualization pipeline.  An important feature of the pipeline architecture is the simplicity of its visual representation. Unfortunately, the notions of an executive, the multiple passes of information flow, flow of metadata (temporal, spatial and contextual) are always depicted in a non-standard fashion in much of software and literature.  While this may not be the right place to discuss a new standard or proposal for this, I think the article would benefit from such a discussion."

$delRange = $d.Range($afterComma, $paraEnd - 1)
$delRange.Delete()

$insPoint = $d.Range($afterComma, $afterComma)
$insPoint.InsertAfter($mergedText)

Write-Host "Done"
